$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new block of 2 rows at row 232 (pushes old row 232.. down by 2) ---
$ws.Rows("232:233").Insert()

# Row 232 (new): Apio, Sin especificar, Primera, date 44846
$ws.Cells.Item(232,1).Value = 11
$ws.Cells.Item(232,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(232,3).Value = "Bíobío"
$ws.Cells.Item(232,4).Value = 44846
$ws.Cells.Item(232,5).Value = 8
$ws.Cells.Item(232,6).Value = 100112017
$ws.Cells.Item(232,7).Value = "Apio"
$ws.Cells.Item(232,8).Value = "Sin especificar"
$ws.Cells.Item(232,9).Value = "Primera"
$ws.Cells.Item(232,10).Value = 150
$ws.Cells.Item(232,11).Value = 8000
$ws.Cells.Item(232,12).Value = 8000
$ws.Cells.Item(232,13).Value = 8000
$ws.Cells.Item(232,14).Value = "$/docena de matas"
$ws.Cells.Item(232,15).Value = "Región de Coquimbo"
$ws.Cells.Item(232,16).Value = 1333
$ws.Cells.Item(232,17).Value = 6
$ws.Cells.Item(232,18).Value = "Hortaliza"

# Row 233 (new): Apio, Sin especificar, Segunda, date 44846
$ws.Cells.Item(233,1).Value = 11
$ws.Cells.Item(233,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(233,3).Value = "Bíobío"
$ws.Cells.Item(233,4).Value = 44846
$ws.Cells.Item(233,5).Value = 8
$ws.Cells.Item(233,6).Value = 100112017
$ws.Cells.Item(233,7).Value = "Apio"
$ws.Cells.Item(233,8).Value = "Sin especificar"
$ws.Cells.Item(233,9).Value = "Segunda"
$ws.Cells.Item(233,10).Value = 180
$ws.Cells.Item(233,11).Value = 6500
$ws.Cells.Item(233,12).Value = 6500
$ws.Cells.Item(233,13).Value = 6500
$ws.Cells.Item(233,14).Value = "$/docena de matas"
$ws.Cells.Item(233,15).Value = "Región de Coquimbo"
$ws.Cells.Item(233,16).Value = 1083
$ws.Cells.Item(233,17).Value = 6
$ws.Cells.Item(233,18).Value = "Hortaliza"

# --- Insert second new block of 2 rows at row 244 (final numbering, after first insert) ---
$ws.Rows("244:245").Insert()

# Row 244 (new): Apio, Americana (o), Primera, date 44845
$ws.Cells.Item(244,1).Value = 11
$ws.Cells.Item(244,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(244,3).Value = "Bíobío"
$ws.Cells.Item(244,4).Value = 44845
$ws.Cells.Item(244,5).Value = 8
$ws.Cells.Item(244,6).Value = 100112017
$ws.Cells.Item(244,7).Value = "Apio"
$ws.Cells.Item(244,8).Value = "Americana (o)"
$ws.Cells.Item(244,9).Value = "Primera"
$ws.Cells.Item(244,10).Value = 270
$ws.Cells.Item(244,11).Value = 7500
$ws.Cells.Item(244,12).Value = 8000
$ws.Cells.Item(244,13).Value = 7722
$ws.Cells.Item(244,14).Value = "$/docena de matas"
$ws.Cells.Item(244,15).Value = "Región de Coquimbo"
$ws.Cells.Item(244,16).Value = 1287
$ws.Cells.Item(244,17).Value = 6
$ws.Cells.Item(244,18).Value = "Hortaliza"

# Row 245 (new): Apio, Americana (o), Segunda, date 44845
$ws.Cells.Item(245,1).Value = 11
$ws.Cells.Item(245,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(245,3).Value = "Bíobío"
$ws.Cells.Item(245,4).Value = 44845
$ws.Cells.Item(245,5).Value = 8
$ws.Cells.Item(245,6).Value = 100112017
$ws.Cells.Item(245,7).Value = "Apio"
$ws.Cells.Item(245,8).Value = "Americana (o)"
$ws.Cells.Item(245,9).Value = "Segunda"
$ws.Cells.Item(245,10).Value = 250
$ws.Cells.Item(245,11).Value = 6500
$ws.Cells.Item(245,12).Value = 6500
$ws.Cells.Item(245,13).Value = 6500
$ws.Cells.Item(245,14).Value = "$/docena de matas"
$ws.Cells.Item(245,15).Value = "Región de Coquimbo"
$ws.Cells.Item(245,16).Value = 1083
$ws.Cells.Item(245,17).Value = 6
$ws.Cells.Item(245,18).Value = "Hortaliza"
